$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1 ("Components")
# ---------------------------------------------------------------------------

# Row 3 (Teensy 3.1): fill in unit price + supplier/URL (was purchased, now has a price)
$ws1.Range("C3").Value = 16.99
$ws1.Range("E3").Value = "CC"
$ws1.Range("H3").Value = "https://www.coolcomponents.co.uk/teensy-3-1.html"

# Row 11: swap boost converter part (TPS63061DSCR -> LMR62421XMF)
$ws1.Range("A11").Value = "LMR62421XMF"
$ws1.Range("C11").Value = 1.95
$ws1.Range("F11").Value = "SOT-23"
$ws1.Range("G11").Value = "Boost converter for motors."
$ws1.Range("H11").Value = "http://uk.farnell.com/texas-instruments/lmr62421xmf/boost-1a-24vout/dp/2064678"

# Insert 3 new rows below row 12 (pushes the totals block down, keeping the blank
# separator row before "Subtotal:")
$ws1.Range("A13:H15").Insert()

# Row 13: JST connector for LiPo
$ws1.Range("A13").Value = "JST side entry 2pin 2mm"
$ws1.Range("B13").Value = 1
$ws1.Range("C13").Value = 0.231
$ws1.Range("D13").Formula = "=B13*C13"
$ws1.Range("E13").Value = "Farnell"
$ws1.Range("F13").Value = "N/A SMD"
$ws1.Range("G13").Value = "JST connector for LiPo"
$ws1.Range("H13").Value = "http://uk.farnell.com/jst-japan-solderless-terminals/s2b-ph-sm4-tb-lf-sn/connector-header-smt-r-a-2mm-2way/dp/9492615?Ntt=S2B-PH-SM4-TB%28LF%29%28SN"

# Row 14: LiPo battery
$ws1.Range("A14").Value = "LiPo"
$ws1.Range("B14").Value = 1
$ws1.Range("C14").Value = 10.99
$ws1.Range("D14").Formula = "=B14*C14"
$ws1.Range("E14").Value = "CC"
$ws1.Range("G14").Value = "Cool Components. From Sparkfun. Self contained discharge protection circuitry "
$ws1.Range("H14").Value = "https://www.coolcomponents.co.uk/lithium-polymer-battery-2000mah.html"

# Row 15: Motors
$ws1.Range("A15").Value = " Motors"
$ws1.Range("B15").Value = 2
$ws1.Range("C15").Value = 3.13
$ws1.Range("D15").Formula = "=B15*C15"
$ws1.Range("E15").Value = "eBay - China"
$ws1.Range("G15").Value = "Consider economics of ordering in larger package but paying VAT and handling charge"
$ws1.Range("H15").Value = "?"

# Column E is a bit wider now that it no longer auto-fits
$ws1.Columns.Item(5).ColumnWidth = 12.1

# Re-pin the shared D formula + totals formulas to the exact ranges (the row
# Insert above nudges them to D2:D16 which is equivalent but not identical)
$ws1.Range("D3:D15").Formula = "=B3*C3"
$ws1.Range("C17").Formula = "=SUM(D2:D15)"
$ws1.Range("C19").Formula = "=0.2*C17"
$ws1.Range("C23").Formula = "=C17+C19+C21"

# Selection / view state
$ws1.Range("G11").Select()

# ---------------------------------------------------------------------------
# Sheet2 ("PCB manufacture" -> "PCB names and values")
# ---------------------------------------------------------------------------

$ws2.Name = "PCB names and values"

$ws2.Columns.Item(1).ColumnWidth = 21.6
$ws2.Columns.Item(2).ColumnWidth = 54.8

$ws2.Range("A1").Value = "Component ID"
$ws2.Range("B1").Value = "Component"
$ws2.Range("C1").Value = "Value"

$ws2.Range("A3").Value = "JP1"
$ws2.Range("B3").Value = "JST connector"

$ws2.Range("A4").Value = "U1"
$ws2.Range("B4").Value = "Teensy 3.1"

$ws2.Range("A5").Value = "IC1"
$ws2.Range("B5").Value = "LMR62421XMF"

$ws2.Range("A6").Select()
$ws2.Activate()
